$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 4620
$ws.Range("B2").Value = 45810
$ws.Range("A3").Value = 4580
$ws.Range("B3").Value = 45810.01041666666
$ws.Range("A4").Value = 4540
$ws.Range("B4").Value = 45810.02083333334
$ws.Range("A5").Value = 4500
$ws.Range("B5").Value = 45810.03125
$ws.Range("A6").Value = 4470
$ws.Range("B6").Value = 45810.04166666666
$ws.Range("A7").Value = 4440
$ws.Range("B7").Value = 45810.05208333334
$ws.Range("A8").Value = 4430
$ws.Range("B8").Value = 45810.0625
$ws.Range("A9").Value = 4430
$ws.Range("B9").Value = 45810.07291666666
$ws.Range("A10").Value = 4430
$ws.Range("B10").Value = 45810.08333333334
$ws.Range("A11").Value = 4440
$ws.Range("B11").Value = 45810.09375
$ws.Range("A12").Value = 4450
$ws.Range("B12").Value = 45810.10416666666
$ws.Range("A13").Value = 4460
$ws.Range("B13").Value = 45810.11458333334
$ws.Range("A14").Value = 4470
$ws.Range("B14").Value = 45810.125
$ws.Range("A15").Value = 4490
$ws.Range("B15").Value = 45810.13541666666
$ws.Range("A16").Value = 4510
$ws.Range("B16").Value = 45810.14583333334
$ws.Range("A17").Value = 4530
$ws.Range("B17").Value = 45810.15625
$ws.Range("A18").Value = 4570
$ws.Range("B18").Value = 45810.16666666666
$ws.Range("A19").Value = 4630
$ws.Range("B19").Value = 45810.17708333334
$ws.Range("A20").Value = 4700
$ws.Range("B20").Value = 45810.1875
$ws.Range("A21").Value = 4800
$ws.Range("B21").Value = 45810.19791666666
$ws.Range("A22").Value = 4910
$ws.Range("B22").Value = 45810.20833333334
$ws.Range("A23").Value = 5030
$ws.Range("B23").Value = 45810.21875
$ws.Range("A24").Value = 5160
$ws.Range("B24").Value = 45810.22916666666
$ws.Range("A25").Value = 5290
$ws.Range("B25").Value = 45810.23958333334
$ws.Range("A26").Value = 5400
$ws.Range("B26").Value = 45810.25
$ws.Range("A27").Value = 5500
$ws.Range("B27").Value = 45810.26041666666
$ws.Range("A28").Value = 5580
$ws.Range("B28").Value = 45810.27083333334
$ws.Range("A29").Value = 5620
$ws.Range("B29").Value = 45810.28125
$ws.Range("A30").Value = 5640
$ws.Range("B30").Value = 45810.29166666666
$ws.Range("A31").Value = 5630
$ws.Range("B31").Value = 45810.30208333334
$ws.Range("A32").Value = 5590
$ws.Range("B32").Value = 45810.3125
$ws.Range("A33").Value = 5530
$ws.Range("B33").Value = 45810.32291666666
$ws.Range("A34").Value = 5460
$ws.Range("B34").Value = 45810.33333333334
$ws.Range("A35").Value = 5380
$ws.Range("B35").Value = 45810.34375
$ws.Range("A36").Value = 5300
$ws.Range("B36").Value = 45810.35416666666
$ws.Range("A37").Value = 5240
$ws.Range("B37").Value = 45810.36458333334
$ws.Range("A38").Value = 5170
$ws.Range("B38").Value = 45810.375
$ws.Range("A39").Value = 5110
$ws.Range("B39").Value = 45810.38541666666
$ws.Range("A40").Value = 5060
$ws.Range("B40").Value = 45810.39583333334
$ws.Range("A41").Value = 5030
$ws.Range("B41").Value = 45810.40625
$ws.Range("A42").Value = 5000
$ws.Range("B42").Value = 45810.41666666666
$ws.Range("A43").Value = 4980
$ws.Range("B43").Value = 45810.42708333334
$ws.Range("A44").Value = 4960
$ws.Range("B44").Value = 45810.4375
$ws.Range("A45").Value = 4950
$ws.Range("B45").Value = 45810.44791666666
$ws.Range("A46").Value = 4950
$ws.Range("B46").Value = 45810.45833333334
$ws.Range("A47").Value = 4950
$ws.Range("B47").Value = 45810.46875
$ws.Range("A48").Value = 4960
$ws.Range("B48").Value = 45810.47916666666
$ws.Range("A49").Value = 4980
$ws.Range("B49").Value = 45810.48958333334
$ws.Range("A50").Value = 5000
$ws.Range("B50").Value = 45810.5
$ws.Range("A51").Value = 5030
$ws.Range("B51").Value = 45810.51041666666
$ws.Range("A52").Value = 5070
$ws.Range("B52").Value = 45810.52083333334
$ws.Range("A53").Value = 5100
$ws.Range("B53").Value = 45810.53125
$ws.Range("A54").Value = 5140
$ws.Range("B54").Value = 45810.54166666666
$ws.Range("A55").Value = 5160
$ws.Range("B55").Value = 45810.55208333334
$ws.Range("A56").Value = 5190
$ws.Range("B56").Value = 45810.5625
$ws.Range("A57").Value = 5210
$ws.Range("B57").Value = 45810.57291666666
$ws.Range("A58").Value = 5230
$ws.Range("B58").Value = 45810.58333333334
$ws.Range("A59").Value = 5250
$ws.Range("B59").Value = 45810.59375
$ws.Range("A60").Value = 5270
$ws.Range("B60").Value = 45810.60416666666
$ws.Range("A61").Value = 5300
$ws.Range("B61").Value = 45810.61458333334
$ws.Range("A62").Value = 5340
$ws.Range("B62").Value = 45810.625
$ws.Range("A63").Value = 5390
$ws.Range("B63").Value = 45810.63541666666
$ws.Range("A64").Value = 5460
$ws.Range("B64").Value = 45810.64583333334
$ws.Range("A65").Value = 5520
$ws.Range("B65").Value = 45810.65625
$ws.Range("A66").Value = 5590
$ws.Range("B66").Value = 45810.66666666666
$ws.Range("A67").Value = 5650
$ws.Range("B67").Value = 45810.67708333334
$ws.Range("A68").Value = 5720
$ws.Range("B68").Value = 45810.6875
$ws.Range("A69").Value = 5790
$ws.Range("B69").Value = 45810.69791666666
$ws.Range("A70").Value = 5850
$ws.Range("B70").Value = 45810.70833333334
$ws.Range("A71").Value = 5920
$ws.Range("B71").Value = 45810.71875
$ws.Range("A72").Value = 6010
$ws.Range("B72").Value = 45810.72916666666
$ws.Range("A73").Value = 6110
$ws.Range("B73").Value = 45810.73958333334
$ws.Range("A74").Value = 6200
$ws.Range("B74").Value = 45810.75
$ws.Range("A75").Value = 6300
$ws.Range("B75").Value = 45810.76041666666
$ws.Range("A76").Value = 6400
$ws.Range("B76").Value = 45810.77083333334
$ws.Range("A77").Value = 6480
$ws.Range("B77").Value = 45810.78125
$ws.Range("A78").Value = 6580
$ws.Range("B78").Value = 45810.79166666666
$ws.Range("A79").Value = 6690
$ws.Range("B79").Value = 45810.80208333334
$ws.Range("A80").Value = 6760
$ws.Range("B80").Value = 45810.8125
$ws.Range("A81").Value = 6840
$ws.Range("B81").Value = 45810.82291666666
$ws.Range("A82").Value = 6860
$ws.Range("B82").Value = 45810.83333333334
$ws.Range("A83").Value = 6840
$ws.Range("B83").Value = 45810.84375
$ws.Range("A84").Value = 6780
$ws.Range("B84").Value = 45810.85416666666
$ws.Range("A85").Value = 6700
$ws.Range("B85").Value = 45810.86458333334
$ws.Range("A86").Value = 6590
$ws.Range("B86").Value = 45810.875
$ws.Range("A87").Value = 6460
$ws.Range("B87").Value = 45810.88541666666
$ws.Range("A88").Value = 6300
$ws.Range("B88").Value = 45810.89583333334
$ws.Range("A89").Value = 6140
$ws.Range("B89").Value = 45810.90625
$ws.Range("A90").Value = 5970
$ws.Range("B90").Value = 45810.91666666666
$ws.Range("A91").Value = 5820
$ws.Range("B91").Value = 45810.92708333334
$ws.Range("A92").Value = 5660
$ws.Range("B92").Value = 45810.9375
$ws.Range("A93").Value = 5520
$ws.Range("B93").Value = 45810.94791666666
$ws.Range("A94").Value = 5350
$ws.Range("B94").Value = 45810.95833333334
$ws.Range("A95").Value = 5300
$ws.Range("B95").Value = 45810.96875
$ws.Range("A96").Value = 5250
$ws.Range("B96").Value = 45810.97916666666
$ws.Range("A97").Value = 5230
$ws.Range("B97").Value = 45810.98958333334
